# Fix set max budget bug
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# City name correction: Bucharest -> Palma de Mallorca (row 3, column A)
$ws.Range("A3").Value = "Palma de Mallorca"

# Max hotel price per night for that city: 1456 -> 3000 (row 3, column E)
$ws.Range("E3").Value = 3000

# Move the active selection to E6 (matches the saved sheet view selection)
$ws.Range("E6").Select()
